# ADD gb.py and btes.py
#
# Each worksheet (one per projection year: 2025, 2030, 2035, 2040, 2045, 2050)
# gains two new technology columns:
#   - "gb"   inserted right after "eb"   (becomes new column B)
#   - "btes" inserted right after "ttes" (becomes new column N, right before "ites")
# and every sheet's row-2 data values are refreshed to the re-run results.

$wb = $excel.ActiveWorkbook

# Row-2 values (A..O) for each sheet, in final post-edit column order:
# eb, gb, hp, st, wi, ieh, chp, ac, ab_ct, ab_hp, cp_ct, cp_hp, ttes, btes, ites
$rowValues = @{
    "2025" = @(0.008591273047792213, 0, 0.5031130441162281, 0, 0.7500791291218792, 0.02935900691239179, 0, 0.06880566916063373, 0, 0, 0, 0, 0, 1490.305292690596, 3.364628064996621)
    "2030" = @(0, 0, 1.742919626181552, 0, 0, 0, 0, 0.1393700898393662, 0, 0, 0, 0, 0, 132.8466098810368, 6.646311143090884)
    "2035" = @(0, 0, 1.72946932639293, 0, 0, 0.4864420908657058, 0, 0.1261817221222202, 0, 0, 0, 0, 0, 14.98775517455852, 7.285821150245502)
    "2040" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 2.302146184846606, 0)
    "2045" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    "2050" = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if (-not $rowValues.ContainsKey($name)) {
        continue
    }

    # Insert the two new columns (shifts everything from that column rightward).
    # "gb" -> column B (index 2), pushing hp..ites one to the right.
    $ws.Columns.Item(2).Insert()
    # "btes" -> column N (index 14), right before "ites" (already shifted to 14).
    $ws.Columns.Item(14).Insert()

    # Header row.
    $ws.Cells.Item(1, 2).Value = "gb"
    $ws.Cells.Item(1, 14).Value = "btes"

    # Data row (row 2), columns A..O (1..15).
    $values = $rowValues[$name]
    for ($col = 1; $col -le 15; $col++) {
        $ws.Cells.Item(2, $col).Value = $values[$col - 1]
    }
}
